$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 text changes from "Empresa Id" to "Entidad Id"
$ws.Range("B1").Value = "Entidad Id"

# Update the selected cell/active cell to B3 (was B14)
$ws.Range("B3").Select()
